# Add two new transaction rows at the top of the data (rows 2-3), shifting the
# existing data rows down by two. Implemented by moving data bottom-up (to avoid
# clobbering not-yet-read cells) instead of a native row Insert, since Insert()
# here would also propagate cell formatting/style (and create empty cell nodes)
# to columns that should remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("E", "N", "P", "T")

# Shift existing data rows 2..8 down to 4..10, working from the bottom up so
# that a source row is always read before it gets overwritten.
for ($srcRow = 8; $srcRow -ge 2; $srcRow--) {
    $dstRow = $srcRow + 2
    foreach ($col in $cols) {
        $srcCell = $ws.Range("$col$srcRow")
        $dstCell = $ws.Range("$col$dstRow")
        $dstCell.Value = $srcCell.Value2
    }
}

# Populate the two newly inserted rows with the new transaction data.
$ws.Range("E2").Value = "Withdrawal"
$ws.Range("N2").Value = "Crypto"
$ws.Range("P2").Value = "USDT ERC 20"
$ws.Range("T2").Value = 1000.0771999999999

$ws.Range("E3").Value = "Withdrawal"
$ws.Range("N3").Value = "Crypto"
$ws.Range("P3").Value = "ETH"
$ws.Range("T3").Value = 540.65

# Update the view: active cell / selection as recorded after the edit, and the
# window horizontal scroll position.
$ws.Range("R19").Select() | Out-Null
$excel.ActiveWindow.Left = -28920 | Out-Null
